# Update metrics values to handle None values in metrics / compute metrics
# for results without optimization (per commit message).
$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Row 2 (adult)
$ws.Range("C2").Value = 2.597419775019497
$ws.Range("D2").Value = 17.75
$ws.Range("E2").Value = 12.93879281845967
$ws.Range("F2").Value = 12.7541939361204
$ws.Range("G2").Value = 1

# Row 3 (compas)
$ws.Range("C3").Value = 1.057060508177166
$ws.Range("D3").Value = 5.161616161616162
$ws.Range("E3").Value = 2.199146899475081
$ws.Range("F3").Value = 2.163708741019674
$ws.Range("G3").Value = 1

# Row 4 (credit)
$ws.Range("C4").Value = 0.1270085289020096
$ws.Range("D4").Value = 9.06
$ws.Range("E4").Value = 0.2420452041110067
$ws.Range("F4").Value = 0.06782578315039454
$ws.Range("G4").Value = 0.1025827849439245

$wb.Save()
